$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revised Day 10 solution run times
$ws.Range("B14").Value = 0.0131570000085048
$ws.Range("C14").Value = 0.0030266999965533601

# Move the active selection to A29
[void]$ws.Range("A29").Select()
